# The edit: justify ("both") paragraphs 6 through 19 (1-based, Word's
# Paragraphs collection is 1-based) which previously had no explicit
# alignment set. Paragraphs 1-5 and 20 already had jc="both" and are
# left untouched. In addition, the "Result:" (paragraph 17) and
# "Conclusion:" (paragraph 19) heading paragraphs are missing the bold
# formatting that the other headings ("Introduction:", "Method:")
# already carry - this edit fixes that inconsistency at the same time.

$d = $word.ActiveDocument

$wdAlignParagraphJustify = 3

for ($i = 6; $i -le 19; $i++) {
    $para = $d.Paragraphs($i)
    $para.Alignment = $wdAlignParagraphJustify
}

$d.Paragraphs(17).Range.Bold = 1
$d.Paragraphs(19).Range.Bold = 1
